# oceny.xlsx - "Add files via upload" edit
#
# Summary of the change (from the OOXML diff):
#  - Arkusz2 (sheet2) becomes the active/selected tab instead of Arkusz1.
#  - Arkusz1's view keeps the same frozen pane (ySplit=3) but the
#    selection/active cell moves from G12 to B6.
#  - Arkusz2's selection changes from a single cell (D14) to the full B1:B17
#    range, with B1 as the active cell.
#  - Arkusz2 column B (rows 1-17) is rewritten from plain numeric IDs
#    (100201..100217) to long hash-like text strings, which also appends 17
#    new entries to the shared-strings table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Arkusz2: replace the numeric IDs in column B with hash strings ------
$hashes = @(
    "d512c8e3dc08fb2d1e8d22d36e6a162a756a7237cdd1fab6904491f8196d9b17",
    "3d3ac11b66f8a9e53116bfa7dfcbb4c21d8fffa4c1709df6caa5f0b20526545c",
    "ce8d28e7c9a6a49c58b0dc074484d1b85516c66c99753ba24e19a32a341638d8",
    "91beb3262326223ac53c506af824d5c261756f6133463fec4087bdb5b9e3a8d3",
    "20b0c9927a562718a1d45dfe529fcd223b1fb1e5dc45e4b1e4b3961676c5d467",
    "b9f0520d80d855b474e806c9d7e362229d66e4f45a75e0fe053ef9767aa7396d",
    "8848110eebbd2747214e193cf91aac11c05918573c753c4b1ac9d40fc74b5fb0",
    "6bfe66a791a8b40d2ff2650d553ecc01060f920bdc8a0749b63a5c0df7aefd7c",
    "8f075116c5fc13c05cf8e91648543ed9630c3ac396ee0667bab2d227c0f1576b",
    "0cd19803993847d06497b256500a558cd2afed8ae70d40a0e93b2f7c4b0682f9",
    "ddbb2f55350fce6fcfb4f06e3874952ddbbaa8a2a04af60fe3c289253aaa07aa",
    "2e0c827843252d195749b4bd1cda2ea8cbc7a8b2ec77d115d790a722730ce395",
    "fd5d85254fe84121251cc29ba390a6aa6a5a320f42c7d321d7c5e184ed0e6a86",
    "c0b3402652b12e59ccc06c2b1c25cdead09c6d46aa4e96c2b8077abc858fb1b6",
    "6532bd55b86fb02cde9815ed74d71c2ad1ccfaa691b5085245f1855857c742db",
    "3f119e8cb45ea54f68720745a91a7a0b9e6af2f263a6465597fd02b278100793",
    "6fa2d129d3ae414f966a4e67e89861a9a1ca5435610ae8dd4e939e19af8c8b23"
)

for ($i = 0; $i -lt $hashes.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 2).Value = $hashes[$i]
}

# --- View/selection state -------------------------------------------------
# Arkusz1 keeps its frozen header rows, only the selected cell changes.
$ws1.Range("B6").Select()

# Arkusz2 gets selected as a range and becomes the active sheet/tab.
$ws2.Range("B1:B17").Select()
$ws2.Activate()
